$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.018.64'
$ws.Range('E2').Value = '  +0.09%  '

$ws.Range('D3').Value = '1.634.36'
$ws.Range('E3').Value = '  -0.39%  '

$ws.Range('E4').Value = '  +0.16%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '214.07'
$ws.Range('D5').NumberFormat = 'General'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -0.87%  '

$ws.Range('E6').Value = '  -0.62%  '

$ws.Range('E7').Value = '  +0.21%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.252'
$ws.Range('D8').NumberFormat = 'General'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -1.95%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.0624'
$ws.Range('D9').NumberFormat = 'General'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  -2.08%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '18.55'
$ws.Range('D10').NumberFormat = 'General'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  -5.28%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0793'
$ws.Range('D11').NumberFormat = 'General'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -0.23%  '

$ws.Range('D12').Value = '1.861.04'
$ws.Range('E12').Value = '  -0.35%  '

$ws.Range('E13').Value = '  -1.57%  '

$ws.Range('D14').Value = '1.637.84'
$ws.Range('E14').Value = '  -1.93%  '

$ws.Range('E15').Value = '  -2.69%  '

$ws.Range('D16').Value = '0.0₃0748'
$ws.Range('E16').Value = '  -2.22%  '

$ws.Range('D17').Value = '26.008.85'
$ws.Range('E17').Value = '  +0.45%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '61.96'
$ws.Range('D18').NumberFormat = 'General'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  -1.58%  '

$ws.Range('E19').Value = '  +0.16%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '191.01'
$ws.Range('D20').NumberFormat = 'General'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  -0.98%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '4.27'
$ws.Range('D21').NumberFormat = 'General'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -2.05%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '9.59'
$ws.Range('D22').NumberFormat = 'General'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  -3.45%  '

$ws.Range('E23').Value = '  -1.72%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '0.133'
$ws.Range('D24').NumberFormat = 'General'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +1.30%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '143.46'
$ws.Range('D25').NumberFormat = 'General'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -0.81%  '

$ws.Range('E26').Value = '  +0.14%  '

$ws.Range('E27').Value = '  -1.83%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '6.76'
$ws.Range('D28').NumberFormat = 'General'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -2.33%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '15.25'
$ws.Range('D29').NumberFormat = 'General'
$ws.Range('D29').Style = 'Normal'
$ws.Range('E29').Value = '  -1.92%  '

$ws.Range('E30').Value = '  -1.45%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.0484'
$ws.Range('D31').NumberFormat = 'General'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -3.08%  '

$ws.Range('E32').Value = '  -2.75%  '

$ws.Range('E33').Value = '  -3.96%  '

$ws.Range('E34').Value = '  -1.38%  '

$ws.Range('E35').Value = '  -2.11%  '

$ws.Range('E36').Value = '  -2.90%  '

$ws.Range('D37').Value = '1.132.99'
$ws.Range('E37').Value = '  -0.10%  '

$ws.Range('E38').Value = '  -0.96%  '

$ws.Range('E39').Value = '  -3.18%  '

$ws.Range('E40').Value = '  -1.42%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '98.89'
$ws.Range('D41').NumberFormat = 'General'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.66%  '

$ws.Range('E42').Value = '  -1.18%  '

$ws.Range('E43').Value = '  -3.63%  '

$ws.Range('D44').Value = '1.771.29'
$ws.Range('E44').Value = '  -0.21%  '

$ws.Range('E45').Value = '  -0.88%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '55.39'
$ws.Range('D46').NumberFormat = 'General'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -2.18%  '

$ws.Range('E47').Value = '  -0.63%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.49'
$ws.Range('D48').NumberFormat = 'General'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  +1.68%  '

$ws.Range('E49').Value = '  -0.30%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '7.53'
$ws.Range('D50').NumberFormat = 'General'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.62%  '

$ws.Range('E51').Value = '  +0.11%  '
